$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity")

# New rows to insert (pairs of A/B values) replacing current rows 13-16 (1-indexed sheet rows)
# Build full target table for rows 2..23 (row 1 header unchanged)

$data = @(
    @("Agriculture, cattling & fishering", "Agriculture, cattling & fishering"),
    @("Chemicals", "Chemicals"),
    @("Electricity by fossil fuels", "Electricity"),
    @("Electricity by nuclear", "Electricity"),
    @("Electricity by other RES", "Electricity"),
    @("Food", "Food"),
    @("Fuels extraction", "Mining & quarrying"),
    @("Fuels refinery", "Fuels refinery"),
    @("Metals", "Metals"),
    @("Mining & quarrying", "Mining & quarrying"),
    @("Other manufacturing", "Other manufacturing"),
    @("Services", "Services"),
    @("Transport", "Transport"),
    @("PV plants", "Other manufacturing"),
    @("PV modules", "Other manufacturing"),
    @("Si-cells", "Other manufacturing"),
    @("Onshore wind plants", "Other manufacturing"),
    @("DFIG generators", "Other manufacturing"),
    @("Offshore wind plants", "Other manufacturing"),
    @("PMG generators", "Other manufacturing"),
    @("Electricity by PV", "Electricity"),
    @("Electricity by wind", "Electricity")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Range("C2").Select()
